$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.00581157207489
$ws.Range("B1").Value = 2.274240493774414
$ws.Range("C1").Value = 4.903514385223389
$ws.Range("D1").Value = 1.649831056594849
$ws.Range("E1").Value = 1.282908916473389
